# Auto-generated Excel COM-interop script to apply the Hyperion_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1019.9091
$ws.Range("I11").Value = 1019.9091
$ws.Range("K11").Value = 1019.9091
$ws.Range("M11").Value = -879.9091
$ws.Range("H12").Value = 2273375.8
$ws.Range("I12").Value = 2841451
$ws.Range("K12").Value = 2841451
$ws.Range("M12").Value = -2841281
$ws.Range("H18").Value = 2540.4443
$ws.Range("I18").Value = 2540.4443
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 2540.4443
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -2256.4443
$ws.Range("N18").ClearContents()
$ws.Range("H28").Value = 1276.9354
$ws.Range("I28").Value = 383.22726
$ws.Range("J28").Value = 3461.5557
$ws.Range("K28").Value = 383.22726
$ws.Range("L28").Value = 3461.5557
$ws.Range("M28").Value = 101.77274
$ws.Range("N28").Value = -4431.5557
$ws.Range("H43").Value = 999
$ws.Range("I43").Value = 999
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 999
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -930
$ws.Range("N43").ClearContents()
$ws.Range("H64").Value = 7235.4375
$ws.Range("I64").Value = 3000
$ws.Range("J64").Value = 7372.0645
$ws.Range("K64").Value = 3000
$ws.Range("L64").Value = 7372.0645
$ws.Range("M64").Value = -2752
$ws.Range("N64").Value = -7868.0645
$ws.Range("H67").Value = 7235.4375
$ws.Range("I67").Value = 3000
$ws.Range("J67").Value = 7372.0645
$ws.Range("K67").Value = 3000
$ws.Range("L67").Value = 7372.0645
$ws.Range("M67").Value = -2142
$ws.Range("N67").Value = -9088.0645
$ws.Range("H69").Value = 6109.5312
$ws.Range("I69").Value = 3991.6667
$ws.Range("J69").Value = 6598.269
$ws.Range("K69").Value = 11975.0001
$ws.Range("L69").Value = 19794.807
$ws.Range("M69").Value = -11101.0001
$ws.Range("N69").Value = -21542.807
$ws.Range("H72").Value = 6109.5312
$ws.Range("I72").Value = 3991.6667
$ws.Range("J72").Value = 6598.269
$ws.Range("K72").Value = 35925.0003
$ws.Range("L72").Value = 59384.421
$ws.Range("M72").Value = -31557.0003
$ws.Range("N72").Value = -68120.421
$ws.Range("H74").Value = 6729.2856
$ws.Range("I74").Value = 5781
$ws.Range("K74").Value = 5781
$ws.Range("M74").Value = -4845
$ws.Range("H76").Value = 6545.2144
$ws.Range("I76").Value = 5520.4287
$ws.Range("K76").Value = 5520.4287
$ws.Range("M76").Value = -5205.4287
$ws.Range("H77").Value = 6729.2856
$ws.Range("I77").Value = 5781
$ws.Range("K77").Value = 28905
$ws.Range("M77").Value = -24225
$ws.Range("H79").Value = 6545.2144
$ws.Range("I79").Value = 5520.4287
$ws.Range("K79").Value = 5520.4287
$ws.Range("M79").Value = -4428.4287
$ws.Range("H80").Value = 5275.1577
$ws.Range("I80").Value = 847.1667
$ws.Range("J80").Value = 7318.846
$ws.Range("K80").Value = 2541.5001
$ws.Range("L80").Value = 21956.538
$ws.Range("M80").Value = -1543.5001
$ws.Range("N80").Value = -23952.538
$ws.Range("H83").Value = 5275.1577
$ws.Range("I83").Value = 847.1667
$ws.Range("J83").Value = 7318.846
$ws.Range("K83").Value = 7624.5003
$ws.Range("L83").Value = 65869.614
$ws.Range("M83").Value = -2632.5003
$ws.Range("N83").Value = -75853.614
$ws.Range("H97").Value = 1325.5
$ws.Range("I97").Value = 802
$ws.Range("J97").Value = 1849
$ws.Range("K97").Value = 2406
$ws.Range("L97").Value = 5547
$ws.Range("M97").Value = -1910
$ws.Range("N97").Value = -6539
$ws.Range("H98").Value = 1268.4348
$ws.Range("I98").Value = 1198.762
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 1198.762
$ws.Range("L98").Value = 2000
$ws.Range("M98").Value = 299.2380000000001
$ws.Range("N98").Value = -4996
$ws.Range("H111").Value = 5292427.5
$ws.Range("I111").Value = 6945608
$ws.Range("K111").Value = 20836824
$ws.Range("M111").Value = -20833757
$ws.Range("H116").Value = 5868.2915
$ws.Range("I116").Value = 3960.5
$ws.Range("J116").Value = 7776.0835
$ws.Range("K116").Value = 3960.5
$ws.Range("L116").Value = 7776.0835
$ws.Range("M116").Value = -518.5
$ws.Range("N116").Value = -14660.0835
$ws.Range("H121").Value = 2233.7693
$ws.Range("J121").Value = 2233.7693
$ws.Range("L121").Value = 6701.3079
$ws.Range("N121").Value = -10195.3079
$ws.Range("H122").Value = 1268.4348
$ws.Range("I122").Value = 1198.762
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 3596.286
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -1146.286
$ws.Range("N122").Value = -10900
$ws.Range("H132").Value = 4206.24
$ws.Range("I132").Value = 4398.087
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 13194.261
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -10664.261
$ws.Range("N132").Value = -11060
$ws.Range("H133").Value = 75000
$ws.Range("J133").Value = 75000
$ws.Range("L133").Value = 75000
$ws.Range("N133").Value = -85120
$ws.Range("H134").Value = 248663.67
$ws.Range("J134").Value = 248663.67
$ws.Range("L134").Value = 248663.67
$ws.Range("N134").Value = -258803.67
$ws.Range("H141").Value = 17435.143
$ws.Range("I141").Value = 23009.2
$ws.Range("J141").Value = 3500
$ws.Range("K141").Value = 69027.6
$ws.Range("L141").Value = 10500
$ws.Range("M141").Value = -63847.60000000001
$ws.Range("N141").Value = -20860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1230552.8
$ws.Range("I2").Value = 1347463.1
$ws.Range("K2").Value = 1347463.1
$ws.Range("M2").Value = -1347350.1
$ws.Range("H45").Value = 6496408.5
$ws.Range("I45").Value = 8929763
$ws.Range("K45").Value = 8929763
$ws.Range("M45").Value = -8929386
$ws.Range("H61").Value = 4840.8335
$ws.Range("I61").Value = 5079.4
$ws.Range("J61").Value = 3648
$ws.Range("K61").Value = 5079.4
$ws.Range("L61").Value = 3648
$ws.Range("M61").Value = -4867.4
$ws.Range("N61").Value = -4072
$ws.Range("H63").Value = 4929.684
$ws.Range("I63").Value = 2342.2727
$ws.Range("K63").Value = 2342.2727
$ws.Range("M63").Value = -1656.2727
$ws.Range("H66").Value = 4929.684
$ws.Range("I66").Value = 2342.2727
$ws.Range("K66").Value = 11711.3635
$ws.Range("M66").Value = -8279.3635
$ws.Range("H74").Value = 25998.87
$ws.Range("I74").Value = 1915.9
$ws.Range("J74").Value = 116310
$ws.Range("K74").Value = 1915.9
$ws.Range("L74").Value = 116310
$ws.Range("M74").Value = -1041.9
$ws.Range("N74").Value = -118058
$ws.Range("H77").Value = 25998.87
$ws.Range("I77").Value = 1915.9
$ws.Range("J77").Value = 116310
$ws.Range("K77").Value = 9579.5
$ws.Range("L77").Value = 581550
$ws.Range("M77").Value = -5211.5
$ws.Range("N77").Value = -590286
$ws.Range("H97").Value = 1678427.5
$ws.Range("I97").Value = 1766739.4
$ws.Range("K97").Value = 1766739.4
$ws.Range("M97").Value = -1766243.4
$ws.Range("H102").Value = 2382611.2
$ws.Range("I102").Value = 2526834
$ws.Range("K102").Value = 2526834
$ws.Range("M102").Value = -2525212
$ws.Range("H116").Value = 1230552.8
$ws.Range("I116").Value = 1347463.1
$ws.Range("K116").Value = 1347463.1
$ws.Range("M116").Value = -1345169.1
$ws.Range("H122").Value = 1226544.9
$ws.Range("I122").Value = 1316940.9
$ws.Range("K122").Value = 3950822.7
$ws.Range("M122").Value = -3948372.7
$ws.Range("H124").Value = 44314.43
$ws.Range("J124").Value = 44314.43
$ws.Range("L124").Value = 44314.43
$ws.Range("N124").Value = -54134.43
$ws.Range("H134").Value = 65800
$ws.Range("J134").Value = 65800
$ws.Range("L134").Value = 65800
$ws.Range("N134").Value = -75940
$ws.Range("H136").Value = 4840.8335
$ws.Range("I136").Value = 5079.4
$ws.Range("J136").Value = 3648
$ws.Range("K136").Value = 15238.2
$ws.Range("L136").Value = 10944
$ws.Range("M136").Value = -12688.2
$ws.Range("N136").Value = -16044
$ws.Range("H140").Value = 82386.5
$ws.Range("J140").Value = 82386.5
$ws.Range("L140").Value = 82386.5
$ws.Range("N140").Value = -92746.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1230552.8
$ws.Range("I3").Value = 1347463.1
$ws.Range("K3").Value = 1347463.1
$ws.Range("M3").Value = -1347349.1
$ws.Range("H36").Value = 1496
$ws.Range("I36").Value = 1338.2858
$ws.Range("K36").Value = 1338.2858
$ws.Range("M36").Value = -804.2858000000001
$ws.Range("H99").Value = 6495232.5
$ws.Range("I99").Value = 7520191
$ws.Range("J99").Value = 3830.6667
$ws.Range("K99").Value = 7520191
$ws.Range("L99").Value = 3830.6667
$ws.Range("M99").Value = -7518693
$ws.Range("N99").Value = -6826.6667
$ws.Range("H105").Value = 12503082
$ws.Range("I105").Value = 12503082
$ws.Range("K105").Value = 12503082
$ws.Range("M105").Value = -12501335
$ws.Range("H132").Value = 85999.4
$ws.Range("J132").Value = 85999.4
$ws.Range("L132").Value = 85999.4
$ws.Range("N132").Value = -96119.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 20000
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H27").Value = 20000
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H31").Value = 14913.127
$ws.Range("I31").Value = 1882.4865
$ws.Range("J31").Value = 26392.5
$ws.Range("K31").Value = 1882.4865
$ws.Range("L31").Value = 26392.5
$ws.Range("M31").Value = -1587.4865
$ws.Range("N31").Value = -26982.5
$ws.Range("H34").Value = 14913.127
$ws.Range("I34").Value = 1882.4865
$ws.Range("J34").Value = 26392.5
$ws.Range("K34").Value = 1882.4865
$ws.Range("L34").Value = 26392.5
$ws.Range("M34").Value = -1680.4865
$ws.Range("N34").Value = -26796.5
$ws.Range("H58").Value = 6742.533
$ws.Range("I58").Value = 8899.25
$ws.Range("K58").Value = 8899.25
$ws.Range("M58").Value = -8696.25
$ws.Range("H62").Value = 3594.2856
$ws.Range("I62").Value = 4037.5
$ws.Range("K62").Value = 4037.5
$ws.Range("M62").Value = -3413.5
$ws.Range("H65").Value = 3594.2856
$ws.Range("I65").Value = 4037.5
$ws.Range("K65").Value = 20187.5
$ws.Range("M65").Value = -17067.5
$ws.Range("H99").Value = 3851.5
$ws.Range("I99").Value = 3851.5
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3851.5
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -2353.5
$ws.Range("N99").ClearContents()
$ws.Range("H105").Value = 793
$ws.Range("I105").Value = 544.8889
$ws.Range("J105").Value = 1239.6
$ws.Range("K105").Value = 544.8889
$ws.Range("L105").Value = 1239.6
$ws.Range("M105").Value = 1202.1111
$ws.Range("N105").Value = -4733.6
$ws.Range("H126").Value = 3851.5
$ws.Range("I126").Value = 3851.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 11554.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -9084.5
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 84599.305
$ws.Range("I132").Value = 54809.844
$ws.Range("K132").Value = 164429.532
$ws.Range("M132").Value = -161899.532
$ws.Range("H134").Value = 1506.8462
$ws.Range("I134").Value = 1407.32
$ws.Range("K134").Value = 4221.96
$ws.Range("M134").Value = -1686.96
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 6742.533
$ws.Range("I136").Value = 8899.25
$ws.Range("K136").Value = 26697.75
$ws.Range("M136").Value = -24147.75
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 39110.914
$ws.Range("J141").Value = 44493.4
$ws.Range("L141").Value = 44493.4
$ws.Range("N141").Value = -54853.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2209.6428
$ws.Range("I3").Value = 1496.7
$ws.Range("J3").Value = 3992
$ws.Range("K3").Value = 4490.1
$ws.Range("L3").Value = 11976
$ws.Range("M3").Value = -4378.1
$ws.Range("N3").Value = -12200
$ws.Range("H4").Value = 52352160
$ws.Range("I4").Value = 66964316
$ws.Range("J4").Value = 14360561
$ws.Range("K4").Value = 200892948
$ws.Range("L4").Value = 43081683
$ws.Range("M4").Value = -200892836
$ws.Range("N4").Value = -43081907
$ws.Range("H59").Value = 4308.3335
$ws.Range("J59").Value = 2872.5
$ws.Range("L59").Value = 8617.5
$ws.Range("N59").Value = -9697.5
$ws.Range("H81").Value = 5806.645
$ws.Range("J81").Value = 6553.778
$ws.Range("L81").Value = 19661.334
$ws.Range("N81").Value = -21907.334
$ws.Range("H84").Value = 5806.645
$ws.Range("J84").Value = 6553.778
$ws.Range("L84").Value = 58984.002
$ws.Range("N84").Value = -70216.00200000001
$ws.Range("H109").Value = 3090
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H113").Value = 4109
$ws.Range("I113").Value = 7909.7144
$ws.Range("J113").Value = 1690.3636
$ws.Range("K113").Value = 23729.1432
$ws.Range("L113").Value = 5071.0908
$ws.Range("M113").Value = -21559.1432
$ws.Range("N113").Value = -9411.0908
$ws.Range("H120").Value = 15504.833
$ws.Range("J120").Value = 13000
$ws.Range("L120").Value = 39000
$ws.Range("N120").Value = -48676
$ws.Range("H132").Value = 2234.1667
$ws.Range("I132").Value = 1800
$ws.Range("J132").Value = 2451.25
$ws.Range("K132").Value = 16200
$ws.Range("L132").Value = 22061.25
$ws.Range("M132").Value = -13670
$ws.Range("N132").Value = -27121.25
$ws.Range("H140").Value = 1321.7059
$ws.Range("I140").Value = 1229.375
$ws.Range("J140").Value = 2799
$ws.Range("K140").Value = 3688.125
$ws.Range("L140").Value = 8397
$ws.Range("M140").Value = 1491.875
$ws.Range("N140").Value = -18757

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 30949.5
$ws.Range("J34").Value = 30949.5
$ws.Range("L34").Value = 30949.5
$ws.Range("N34").Value = -31485.5
$ws.Range("H70").Value = 10535393
$ws.Range("I70").Value = 15387944
$ws.Range("K70").Value = 15387944
$ws.Range("M70").Value = -15387674
$ws.Range("H73").Value = 10535393
$ws.Range("I73").Value = 15387944
$ws.Range("K73").Value = 15387944
$ws.Range("M73").Value = -15387008
$ws.Range("H76").Value = 30949.5
$ws.Range("J76").Value = 30949.5
$ws.Range("L76").Value = 30949.5
$ws.Range("N76").Value = -31579.5
$ws.Range("H79").Value = 30949.5
$ws.Range("J79").Value = 30949.5
$ws.Range("L79").Value = 30949.5
$ws.Range("N79").Value = -33133.5
$ws.Range("H97").Value = 1191192.2
$ws.Range("I97").Value = 1323396.9
$ws.Range("K97").Value = 1323396.9
$ws.Range("M97").Value = -1322900.9
$ws.Range("H99").Value = 5368.5713
$ws.Range("I99").Value = 6155.1665
$ws.Range("J99").Value = 649
$ws.Range("K99").Value = 6155.1665
$ws.Range("L99").Value = 649
$ws.Range("M99").Value = -3909.1665
$ws.Range("N99").Value = -5141
$ws.Range("H117").Value = 49993
$ws.Range("J117").Value = 49993
$ws.Range("L117").Value = 49993
$ws.Range("N117").Value = -56877
$ws.Range("H126").Value = 3222857.8
$ws.Range("I126").Value = 2068414.4
$ws.Range("K126").Value = 6205243.199999999
$ws.Range("M126").Value = -6202773.199999999
$ws.Range("H132").Value = 2142.0327
$ws.Range("I132").Value = 1914.25
$ws.Range("J132").Value = 3458.111
$ws.Range("K132").Value = 5742.75
$ws.Range("L132").Value = 10374.333
$ws.Range("M132").Value = -3212.75
$ws.Range("N132").Value = -15434.333
$ws.Range("H133").Value = 55000
$ws.Range("J133").Value = 55000
$ws.Range("L133").Value = 55000
$ws.Range("N133").Value = -65120
$ws.Range("H136").Value = 14975.321
$ws.Range("J136").Value = 14975.321
$ws.Range("L136").Value = 44925.963
$ws.Range("N136").Value = -50025.963
$ws.Range("H140").Value = 83362.73
$ws.Range("J140").Value = 83362.73
$ws.Range("L140").Value = 83362.73
$ws.Range("N140").Value = -93722.73
$ws.Range("H141").Value = 53856.5
$ws.Range("J141").Value = 71170.4
$ws.Range("L141").Value = 71170.4
$ws.Range("N141").Value = -81530.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3284.175
$ws.Range("I46").Value = 1652.3334
$ws.Range("J46").Value = 3572.147
$ws.Range("K46").Value = 1652.3334
$ws.Range("L46").Value = 3572.147
$ws.Range("M46").Value = -1464.3334
$ws.Range("N46").Value = -3948.147
$ws.Range("H55").Value = 2014.3846
$ws.Range("I55").Value = 1527.8334
$ws.Range("J55").Value = 2431.4285
$ws.Range("K55").Value = 1527.8334
$ws.Range("L55").Value = 2431.4285
$ws.Range("M55").Value = -1354.8334
$ws.Range("N55").Value = -2777.4285
$ws.Range("H68").Value = 4403.5557
$ws.Range("I68").Value = 3967.3333
$ws.Range("K68").Value = 3967.3333
$ws.Range("M68").Value = -3218.3333
$ws.Range("H71").Value = 4403.5557
$ws.Range("I71").Value = 3967.3333
$ws.Range("K71").Value = 19836.6665
$ws.Range("M71").Value = -16092.6665
$ws.Range("H93").Value = 13342919
$ws.Range("I93").Value = 18519970
$ws.Range("J93").Value = 30500.572
$ws.Range("K93").Value = 18519970
$ws.Range("L93").Value = 30500.572
$ws.Range("M93").Value = -18518722
$ws.Range("N93").Value = -32996.572
$ws.Range("H106").Value = 14287.5
$ws.Range("J106").Value = 14287.5
$ws.Range("L106").Value = 14287.5
$ws.Range("N106").Value = -16811.5
$ws.Range("H122").Value = 5337.5386
$ws.Range("I122").Value = 3898.3333
$ws.Range("J122").Value = 6099.4707
$ws.Range("K122").Value = 11694.9999
$ws.Range("L122").Value = 18298.4121
$ws.Range("M122").Value = -9244.999899999999
$ws.Range("N122").Value = -23198.4121
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H132").Value = 6385.283
$ws.Range("I132").Value = 6462.523
$ws.Range("K132").Value = 19387.569
$ws.Range("M132").Value = -16857.569
$ws.Range("H133").Value = 118527.2
$ws.Range("J133").Value = 118527.2
$ws.Range("L133").Value = 118527.2
$ws.Range("N133").Value = -123587.2
$ws.Range("H139").Value = 75957
$ws.Range("J139").Value = 75957
$ws.Range("L139").Value = 75957
$ws.Range("N139").Value = -86237
$ws.Range("H140").Value = 87723
$ws.Range("J140").Value = 120923.5
$ws.Range("L140").Value = 120923.5
$ws.Range("N140").Value = -131283.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2192.5
$ws.Range("I96").Value = 2080.2856
$ws.Range("J96").Value = 2323.4167
$ws.Range("K96").Value = 2080.2856
$ws.Range("L96").Value = 2323.4167
$ws.Range("M96").Value = -707.2856000000002
$ws.Range("N96").Value = -5069.4167
$ws.Range("H107").Value = 47620744
$ws.Range("I107").Value = 111111860
$ws.Range("J107").Value = 2413.5833
$ws.Range("K107").Value = 333335580
$ws.Range("L107").Value = 7240.749899999999
$ws.Range("M107").Value = -333333660
$ws.Range("N107").Value = -11080.7499
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H122").Value = 3435.862
$ws.Range("I122").Value = 2101.15
$ws.Range("K122").Value = 6303.450000000001
$ws.Range("M122").Value = -3853.450000000001
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H132").Value = 18724068
$ws.Range("I132").Value = 25645662
$ws.Range("J132").Value = 727927.2
$ws.Range("K132").Value = 76936986
$ws.Range("L132").Value = 2183781.6
$ws.Range("M132").Value = -76934456
$ws.Range("N132").Value = -2188841.6
$ws.Range("H136").Value = 1032.5
$ws.Range("J136").Value = 2261.889
$ws.Range("L136").Value = 6785.667
$ws.Range("N136").Value = -11885.667
$ws.Range("H137").Value = 86333
$ws.Range("J137").Value = 86333
$ws.Range("L137").Value = 86333
$ws.Range("N137").Value = -96533
$ws.Range("H138").Value = 89606.75
$ws.Range("J138").Value = 89606.75
$ws.Range("L138").Value = 89606.75
$ws.Range("N138").Value = -99886.75
